# Leave Card update: extend the attendance table (Table1) with 17 new
# periods (rows 87-103), mirroring the formatting of the last existing
# row (86) and filling in the two leave entries that were recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# 1. Grow the table by 17 rows (A8:K86 -> A8:K103). ListRows.Add keeps the
#    table/autofilter ref, header and totals metadata in sync.
for ($i = 0; $i -lt 17; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# 2. Clone the formatting (number formats, borders, alignment) of the prior
#    last row onto every new row in one shot.
$ws.Range("A86:K86").Copy()
$ws.Range("A87:K103").PasteSpecial(-4122)

# 3. The K column ("REMARKS" date) only carries its date-number-format on
#    rows that actually hold a remark date (88-90 below); everywhere else
#    it should keep the same blank style the rest of the table uses.
$ws.Range("K9").Copy()
$ws.Range("K87").PasteSpecial(-4122)
$ws.Range("K91:K103").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 4. PERIOD column (A) - one entry per pay period.
$ws.Range("A87").Value = 45170
$ws.Range("A88").Value = 45200
$ws.Range("A91").Value = 45231
$ws.Range("A92").Value = 45261
$ws.Range("A93").Value = 45292
$ws.Range("A94").Value = 45323
$ws.Range("A95").Value = 45352
$ws.Range("A96").Value = 45383
$ws.Range("A97").Value = 45413
$ws.Range("A98").Value = 45444
$ws.Range("A99").Value = 45474
$ws.Range("A100").Value = 45505
$ws.Range("A101").Value = 45536
$ws.Range("A102").Value = 45566
$ws.Range("A103").Value = 45597

# 5. PARTICULARS column (B) - leave remarks referencing existing shared text.
$ws.Range("B88").Value = "SP(1-0-0)"
$ws.Range("B89").Value = "SL(1-0-0)"
$ws.Range("B90").Value = "SL(1-0-0)"

# 6. EARNED column (C) - credits earned for the two new SP/SL entries.
$ws.Range("C87").Value = 1.25
$ws.Range("C88").Value = 1.25

# 7. Absence Undertime W/O Pay (H) for the sick-leave rows.
$ws.Range("H89").Value = 1
$ws.Range("H90").Value = 1

# 8. REMARKS date column (K) for the three leave rows.
$ws.Range("K88").Value = 45221
$ws.Range("K89").Value = 45230
$ws.Range("K90").Value = 45222

# 9. "EARNED " calculated column (G) mirrors EARNED via the table formula -
#    restore the formula on every new row (PasteSpecial above only copied
#    formats, not formulas). Set cell-by-cell so Excel doesn't collapse
#    them into one shared-formula block.
for ($r = 87; $r -le 103; $r++) {
    $ws.Range("G$r").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
}
